{"js": "// Update the worksheet date and regenerate every arithmetic problem/answer\n// pair in the table. Each entry below is [oldText, newText] taken from the\n// canonical OOXML diff; every oldText string is unique in the document, so a\n// plain exact-text search+replace is unambiguous for every cell.\nconst replacements = [\n  [\"2024-01-09 Tuesday\", \"2024-01-10 Wednesday\"],\n  [\"5+54=59\", \"67-8=59\"],\n  [\"7+67=74\", \"33-0=33\"],\n  [\"55-27=28\", \"93-40=53\"],\n  [\"61-16=45\", \"42+17=59\"],\n  [\"21+28=49\", \"35+38=73\"],\n  [\"80-52=28\", \"76-33=43\"],\n  [\"68-49=19\", \"51-24=27\"],\n  [\"17+58=75\", \"53-5=48\"],\n  [\"80+19=99\", \"98-7=91\"],\n  [\"11-4=7\", \"45-12=33\"],\n  [\"94-7=87\", \"72-67=5\"],\n  [\"3+61=64\", \"96-5=91\"],\n  [\"87-22=65\", \"64-9=55\"],\n  [\"34+23=57\", \"10+30=40\"],\n  [\"84-11=73\", \"76-49=27\"],\n  [\"70-44=26\", \"42+35=77\"],\n  [\"53+32=85\", \"6+57=63\"],\n  [\"1+10=11\", \"90-8=82\"],\n  [\"78-22=56\", \"35+44=79\"],\n  [\"87-48=39\", \"58+32=90\"],\n  [\"50-13=37\", \"2+30=32\"],\n  [\"26+59=85\", \"29+29=58\"],\n  [\"59+14=73\", \"71+14=85\"],\n  [\"98-29=69\", \"99-27=72\"],\n  [\"50-33=17\", \"2+2=4\"],\n  [\"89+3=92\", \"75-27=48\"],\n  [\"33-15=18\", \"91-84=7\"],\n  [\"80-1=79\", \"43-7=36\"],\n  [\"75-14=61\", \"9+2=11\"],\n  [\"49-2=47\", \"85-69=16\"],\n  [\"64-56=8\", \"79-39=40\"],\n  [\"90-1=89\", \"5+74=79\"],\n  [\"37+41=78\", \"80+13=93\"],\n  [\"2+44=46\", \"72+23=95\"],\n  [\"17+79=96\", \"45+42=87\"],\n  [\"22-17=5\", \"81-36=45\"],\n  [\"39+29=68\", \"98-45=53\"],\n  [\"74-44=30\", \"15+60=75\"],\n  [\"0+71=71\", \"51+26=77\"],\n  [\"26+18=44\", \"96-40=56\"],\n  [\"4+88=92\", \"0+34=34\"],\n  [\"40+0=40\", \"8+75=83\"],\n  [\"83-1=82\", \"29+48=77\"],\n  [\"12+34=46\", \"60+2=62\"],\n  [\"75+13=88\", \"73-70=3\"],\n  [\"69-42=27\", \"13+7=20\"],\n  [\"77-59=18\", \"91-1=90\"],\n  [\"9+70=79\", \"6+49=55\"],\n  [\"22+67=89\", \"21+11=32\"],\n  [\"20+37=57\", \"24+31=55\"],\n  [\"2+25=27\", \"44+23=67\"],\n  [\"9+58=67\", \"61-1=60\"],\n  [\"1+20=21\", \"45+50=95\"],\n  [\"30+30=60\", \"81-27=54\"],\n  [\"73-56=17\", \"50-46=4\"],\n  [\"67+19=86\", \"86-26=60\"],\n  [\"13+46=59\", \"17+36=53\"],\n  [\"94-88=6\", \"70+2=72\"],\n  [\"18+63=81\", \"79-38=41\"],\n  [\"3+7=10\", \"31+62=93\"],\n  [\"81-24=57\", \"23+2=25\"],\n  [\"95-27=68\", \"69-27=42\"],\n  [\"9-9=0\", \"80-20=60\"],\n  [\"71-6=65\", \"67-56=11\"],\n  [\"71-3=68\", \"76-51=25\"],\n  [\"6+65=71\", \"16+33=49\"],\n  [\"67-49=18\", \"91-51=40\"],\n  [\"47-23=24\", \"29+24=53\"],\n  [\"98-23=75\", \"62+6=68\"],\n  [\"4+6=10\", \"84-36=48\"],\n  [\"36+12=48\", \"2+46=48\"],\n  [\"21-8=13\", \"10+76=86\"],\n  [\"12+73=85\", \"63-2=61\"],\n  [\"5+66=71\", \"53-20=33\"],\n  [\"52+4=56\", \"3+60=63\"],\n  [\"59+11=70\", \"89-7=82\"],\n  [\"94-86=8\", \"32+38=70\"],\n  [\"1+61=62\", \"37+51=88\"],\n  [\"24+68=92\", \"95-45=50\"],\n  [\"6+19=25\", \"18-11=7\"],\n  [\"62-8=54\", \"11+32=43\"],\n  [\"49+26=75\", \"31+38=69\"],\n  [\"64+30=94\", \"36+9=45\"],\n  [\"24+6=30\", \"52-24=28\"],\n  [\"60-49=11\", \"89-70=19\"],\n  [\"58-48=10\", \"58-31=27\"],\n  [\"88-83=5\", \"30-11=19\"],\n  [\"90-56=34\", \"21-7=14\"],\n  [\"6+60=66\", \"66+8=74\"],\n  [\"69+19=88\", \"49-31=18\"],\n  [\"77-42=35\", \"91+2=93\"],\n  [\"8+62=70\", \"49-13=36\"],\n  [\"69-24=45\", \"8+16=24\"],\n  [\"15+5=20\", \"34+14=48\"],\n  [\"27-20=7\", \"69-43=26\"],\n  [\"68-57=11\", \"73+21=94\"],\n  [\"88-62=26\", \"77-49=28\"],\n  [\"17+73=90\", \"95-83=12\"],\n  [\"10+14=24\", \"18-12=6\"],\n  [\"96-51=45\", \"60-55=5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // Exact, case-sensitive, non-wildcard match on the literal old text.\n  const found = body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n    matchWildcards: false,\n  });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and regenerate every arithmetic problem/answer\n# pair in the table. Each entry below is (oldText, newText) taken from the\n# canonical OOXML diff; every oldText string is unique in the document, so a\n# plain exact-text Find/Replace is unambiguous for every cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-09 Tuesday\", \"2024-01-10 Wednesday\"),\n    @(\"5+54=59\", \"67-8=59\"),\n    @(\"7+67=74\", \"33-0=33\"),\n    @(\"55-27=28\", \"93-40=53\"),\n    @(\"61-16=45\", \"42+17=59\"),\n    @(\"21+28=49\", \"35+38=73\"),\n    @(\"80-52=28\", \"76-33=43\"),\n    @(\"68-49=19\", \"51-24=27\"),\n    @(\"17+58=75\", \"53-5=48\"),\n    @(\"80+19=99\", \"98-7=91\"),\n    @(\"11-4=7\", \"45-12=33\"),\n    @(\"94-7=87\", \"72-67=5\"),\n    @(\"3+61=64\", \"96-5=91\"),\n    @(\"87-22=65\", \"64-9=55\"),\n    @(\"34+23=57\", \"10+30=40\"),\n    @(\"84-11=73\", \"76-49=27\"),\n    @(\"70-44=26\", \"42+35=77\"),\n    @(\"53+32=85\", \"6+57=63\"),\n    @(\"1+10=11\", \"90-8=82\"),\n    @(\"78-22=56\", \"35+44=79\"),\n    @(\"87-48=39\", \"58+32=90\"),\n    @(\"50-13=37\", \"2+30=32\"),\n    @(\"26+59=85\", \"29+29=58\"),\n    @(\"59+14=73\", \"71+14=85\"),\n    @(\"98-29=69\", \"99-27=72\"),\n    @(\"50-33=17\", \"2+2=4\"),\n    @(\"89+3=92\", \"75-27=48\"),\n    @(\"33-15=18\", \"91-84=7\"),\n    @(\"80-1=79\", \"43-7=36\"),\n    @(\"75-14=61\", \"9+2=11\"),\n    @(\"49-2=47\", \"85-69=16\"),\n    @(\"64-56=8\", \"79-39=40\"),\n    @(\"90-1=89\", \"5+74=79\"),\n    @(\"37+41=78\", \"80+13=93\"),\n    @(\"2+44=46\", \"72+23=95\"),\n    @(\"17+79=96\", \"45+42=87\"),\n    @(\"22-17=5\", \"81-36=45\"),\n    @(\"39+29=68\", \"98-45=53\"),\n    @(\"74-44=30\", \"15+60=75\"),\n    @(\"0+71=71\", \"51+26=77\"),\n    @(\"26+18=44\", \"96-40=56\"),\n    @(\"4+88=92\", \"0+34=34\"),\n    @(\"40+0=40\", \"8+75=83\"),\n    @(\"83-1=82\", \"29+48=77\"),\n    @(\"12+34=46\", \"60+2=62\"),\n    @(\"75+13=88\", \"73-70=3\"),\n    @(\"69-42=27\", \"13+7=20\"),\n    @(\"77-59=18\", \"91-1=90\"),\n    @(\"9+70=79\", \"6+49=55\"),\n    @(\"22+67=89\", \"21+11=32\"),\n    @(\"20+37=57\", \"24+31=55\"),\n    @(\"2+25=27\", \"44+23=67\"),\n    @(\"9+58=67\", \"61-1=60\"),\n    @(\"1+20=21\", \"45+50=95\"),\n    @(\"30+30=60\", \"81-27=54\"),\n    @(\"73-56=17\", \"50-46=4\"),\n    @(\"67+19=86\", \"86-26=60\"),\n    @(\"13+46=59\", \"17+36=53\"),\n    @(\"94-88=6\", \"70+2=72\"),\n    @(\"18+63=81\", \"79-38=41\"),\n    @(\"3+7=10\", \"31+62=93\"),\n    @(\"81-24=57\", \"23+2=25\"),\n    @(\"95-27=68\", \"69-27=42\"),\n    @(\"9-9=0\", \"80-20=60\"),\n    @(\"71-6=65\", \"67-56=11\"),\n    @(\"71-3=68\", \"76-51=25\"),\n    @(\"6+65=71\", \"16+33=49\"),\n    @(\"67-49=18\", \"91-51=40\"),\n    @(\"47-23=24\", \"29+24=53\"),\n    @(\"98-23=75\", \"62+6=68\"),\n    @(\"4+6=10\", \"84-36=48\"),\n    @(\"36+12=48\", \"2+46=48\"),\n    @(\"21-8=13\", \"10+76=86\"),\n    @(\"12+73=85\", \"63-2=61\"),\n    @(\"5+66=71\", \"53-20=33\"),\n    @(\"52+4=56\", \"3+60=63\"),\n    @(\"59+11=70\", \"89-7=82\"),\n    @(\"94-86=8\", \"32+38=70\"),\n    @(\"1+61=62\", \"37+51=88\"),\n    @(\"24+68=92\", \"95-45=50\"),\n    @(\"6+19=25\", \"18-11=7\"),\n    @(\"62-8=54\", \"11+32=43\"),\n    @(\"49+26=75\", \"31+38=69\"),\n    @(\"64+30=94\", \"36+9=45\"),\n    @(\"24+6=30\", \"52-24=28\"),\n    @(\"60-49=11\", \"89-70=19\"),\n    @(\"58-48=10\", \"58-31=27\"),\n    @(\"88-83=5\", \"30-11=19\"),\n    @(\"90-56=34\", \"21-7=14\"),\n    @(\"6+60=66\", \"66+8=74\"),\n    @(\"69+19=88\", \"49-31=18\"),\n    @(\"77-42=35\", \"91+2=93\"),\n    @(\"8+62=70\", \"49-13=36\"),\n    @(\"69-24=45\", \"8+16=24\"),\n    @(\"15+5=20\", \"34+14=48\"),\n    @(\"27-20=7\", \"69-43=26\"),\n    @(\"68-57=11\", \"73+21=94\"),\n    @(\"88-62=26\", \"77-49=28\"),\n    @(\"17+73=90\", \"95-83=12\"),\n    @(\"10+14=24\", \"18-12=6\"),\n    @(\"96-51=45\", \"60-55=5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # MatchCase:$true, MatchWholeWord:$false, MatchWildcards:$false,\n    # Wrap:wdFindContinue(1), Format:$false, ReplaceWith:$newText,\n    # Replace:wdReplaceAll(2)\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
